$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: chandu@gmail.com / keshav12345 (with mailto hyperlink on A3)
$ws.Range("A3").Value = "chandu@gmail.com"
$ws.Range("B3").Value = "keshav12345"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:chandu@gmail.com") | Out-Null

# New row 4: supriya@gmail.com / supriya12345 (with mailto hyperlink on A4)
$ws.Range("A4").Value = "supriya@gmail.com"
$ws.Range("B4").Value = "supriya12345"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:supriya@gmail.com") | Out-Null

# Match the author's final selection (active cell B4)
$ws.Range("B4").Select() | Out-Null
